$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.163.91'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '1.825.27'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("D4").Value = "'0.9980"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'234.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.85%  '
$ws.Range("D6").Value = "'0.6002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.14%  '
$ws.Range("D7").Value = "'0.9999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -4.95%  '
$ws.Range("D9").Value = "'0.2792"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'23.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.63%  '
$ws.Range("D11").Value = "'0.07602"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.43%  '
$ws.Range("D12").Value = '1.828.25'
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("D13").Value = "'4.790"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.26%  '
$ws.Range("D14").Value = "'0.6287"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.76%  '
$ws.Range("D15").Value = "'0.000009921"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.95%  '
$ws.Range("D16").Value = '2.079.54'
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").Value = "'78.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.43%  '
$ws.Range("D18").Value = "'5.847"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.91%  '
$ws.Range("D19").Value = '29.163.65'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = "'226.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("D21").Value = "'1.0000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D23").Value = "'6.985"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.17%  '
$ws.Range("D24").Value = "'0.9994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").Value = "'155.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'8.022"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.21%  '
$ws.Range("E27").Value = '  -3.44%  '
$ws.Range("D28").Value = "'16.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("D29").Value = "'1.489"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.23%  '
$ws.Range("D30").Value = "'0.06229"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -14.47%  '
$ws.Range("D31").Value = "'1.446"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.02%  '
$ws.Range("E32").Value = '  -5.20%  '
$ws.Range("D33").Value = "'3.807"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.77%  '
$ws.Range("D34").Value = "'1.121"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.57%  '
$ws.Range("D35").Value = "'1.741"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.22%  '
$ws.Range("D36").Value = "'0.6405"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.94%  '
$ws.Range("D37").Value = "'2.532"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("D38").Value = '1.216.56'
$ws.Range("E38").Value = '  -1.31%  '
$ws.Range("D39").Value = "'2.730"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.91%  '
$ws.Range("D40").Value = "'0.01733"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'6.505"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.95%  '
$ws.Range("D42").Value = "'0.9053"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.38%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").Value = '1.993.23'
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("D45").Value = "'100.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("D46").Value = "'62.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000116"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.02%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = "'1.595"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.90%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'8.473"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.88%  '
$ws.Range("D50").Value = "'0.4544"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("D51").Value = "'0.05497"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.72%  '
